$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking prices/volume refresh for Thu May  4 23:11:34 UTC 2023.
# Column D/E cells are stored as literal text in the workbook (prices like
# "28.873.55" use dot thousand-separators, not valid Excel numbers; volume
# cells are padded percentage strings). Force Text format cell-by-cell before
# each write so Excel keeps the exact literal string instead of silently
# parsing it into a number (which would also truncate trailing zeros such as
# "0.07860" -> 0.0786).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.870.40'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.19%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.878.67'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.95%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.45%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.49'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.23%  '

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.41%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4617'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.10%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3879'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07860'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.56%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9838'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.18%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.78'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.19%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.810.57'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -5.77%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.995'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.39%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.662'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.92%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06985'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.24%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.42'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.31%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.36%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009945'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.81%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.95'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.55%  '

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.19%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.871.69'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.19%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.263'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.81%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.99'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.79%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.103'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.94%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.31'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.40%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.37'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.08%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.920'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.32%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.79'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.71%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -6.78%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09360'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.55%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.9009'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -5.19%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.269'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.23%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.319'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.80%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.247'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.06%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05739'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.63%  '

# Row 36
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.167'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.12%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02075'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.84%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.001'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.37%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.633'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.59%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5657'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.33%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1773'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.690'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.15%  '

# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.247'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.78%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.91'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.95%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5337'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.53%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07042'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.65%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.843'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.16%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.533'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.02%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.45'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.07%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.068'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -5.81%  '

# Row 51
$ws.Range("B51").Value = 'Chiliz'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GSCt2y6YSgO26+chiliz-chz'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1301'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.13%  '
